# Insert a new data row at row 113, pushing the existing rows 113-238 down
# to 114-239 (row 238's data ends up at row 239, growing the used range from
# A1:T238 to A1:T239).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with its data.
$ws.Range("A113").Value = 4
$ws.Range("B113").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C113").Value = "Los Lagos"
$ws.Range("D113").Value = (Get-Date -Year 2022 -Month 8 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E113").Value = 10
$ws.Range("F113").Value = "Fruta"
$ws.Range("G113").Value = 100108
$ws.Range("H113").Value = "Tropicales y subtropicales"
$ws.Range("I113").Value = 100108002
$ws.Range("J113").Value = "Mango"
$ws.Range("K113").Value = "Sin especificar"
$ws.Range("L113").Value = "Primera"
$ws.Range("M113").Value = 200
$ws.Range("N113").Value = 12000
$ws.Range("O113").Value = 13000
$ws.Range("P113").Value = 12500
$ws.Range("Q113").Value = "$/bandeja 4 kilos"
$ws.Range("R113").Value = "Brasil"
$ws.Range("S113").Value = 3125
$ws.Range("T113").Value = 4
